$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.433.30"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.606.83"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.90"
$ws.Range("E5").Value = "  +3.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.26"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.614.08"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.157"
$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("E13").Value = "  +3.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.058.67"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.52"
$ws.Range("E15").Value = "  +5.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.422.91"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +2.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.607.17"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.57"
$ws.Range("E19").Value = "  +10.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.64"
$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.07"
$ws.Range("E21").Value = "  +2.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.91"
$ws.Range("E22").Value = "  +4.67%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.523"
$ws.Range("E24").Value = "  +5.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.24"
$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  +8.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0798"
$ws.Range("E29").Value = "  +3.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").Value = "  +10.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.40"
$ws.Range("E31").Value = "  +4.34%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.46"
$ws.Range("E33").Value = "  +3.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.47"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  +6.09%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("E36").Value = "  +9.44%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.988"
$ws.Range("E37").Value = "  +8.67%  "

$ws.Range("E38").Value = "  +9.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.14"
$ws.Range("E39").Value = "  +1.47%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "314.59"
$ws.Range("E40").Value = "  +7.74%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  +6.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.842"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.02"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.08"
$ws.Range("E45").Value = "  +6.25%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.87"
$ws.Range("E47").Value = "  +4.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.608"
$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0551"
$ws.Range("E49").Value = "  +2.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.14"
$ws.Range("E50").Value = "  +7.19%  "

$ws.Range("E51").Value = "  +2.18%  "
